$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fromCSV")

# Helper: write a value into a cell as TEXT (shared string) even when the
# text looks like a number, without disturbing the cell's existing style.
# A direct `$ws.Range(...).Value = "16"` would be auto-coerced to a Number
# by Excel's normal type inference. Building the text with a formula and
# pasting only the *values* back in keeps the original cell style intact
# while still landing as a text/string cell.
function Set-TextValue {
    param($cellRef, [string]$text)
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

# 1) short-url column (B) for every data row: "fr9oS0" -> "d7GD8n"
for ($r = 2; $r -le 319; $r++) {
    $ws.Cells.Item($r, 2).Value = "d7GD8n"
}

# 2) Count corrections stored as text in the sheet
Set-TextValue "O306" "16"

Set-TextValue "O310" "7"

Set-TextValue "N311" "427"
Set-TextValue "O311" "66"

Set-TextValue "N313" "190"
Set-TextValue "O313" "164"

Set-TextValue "N315" "17"
Set-TextValue "O315" "5"

Set-TextValue "S316" "467"

Set-TextValue "T317" "0"

Set-TextValue "O318" "8"
